$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data
$ws.Range("A6").Value = "Nitha"
$ws.Range("B6").Value = "Sonith Ajmal"
$ws.Range("C6").Value = "Samsung"

# Adjust column B width (closest achievable value to the target stored
# width of 11.44140625 characters, since width is stored in pixel-quantized
# steps internally)
$ws.Columns.Item(2).ColumnWidth = 10.73

# Update selection to C6
$ws.Range("C6").Select()
